$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Swap the "Periodo Mora" (E) and "Valor Mora" (F) values between rows 16 and 17.
$ws.Range("E16").Value = "1712"
$ws.Range("F16").Value = 29520

$ws.Range("E17").Value = "1711"
$ws.Range("F17").Value = 10824
